$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'art class clothes target'
$ws.Range("A2").Value = 'art class leggings'
$ws.Range("A3").Value = 'art class pants'
$ws.Range("A4").Value = 'art of smart football'
$ws.Range("A5").Value = 'art panties'
$ws.Range("A6").Value = 'art tights'
$ws.Range("A7").Value = 'arthritis knee support for women'
$ws.Range("A8").Value = 'artritis brace'
$ws.Range("A9").Value = 'asics leggings'
$ws.Range("A10").Value = 'asics leggings women'
$ws.Range("A11").Value = 'athleta black pants'
$ws.Range("A12").Value = 'athleta brand women'
$ws.Range("A13").Value = 'athleta capris for women'
$ws.Range("A14").Value = 'athleta clothing'
$ws.Range("A15").Value = 'athleta clothing women'
$ws.Range("A16").Value = 'athleta clothing women brand'
$ws.Range("A17").Value = 'athleta clothing women pants'
$ws.Range("A18").Value = 'athleta clothing women tops'
$ws.Range("A19").Value = 'athleta dresses'
$ws.Range("A20").Value = 'athleta leggings'
$ws.Range("A21").Value = 'athleta leggings capri'
$ws.Range("A22").Value = 'athleta leggings for women'
$ws.Range("A23").Value = 'athleta leggings for women with pockets'
$ws.Range("A24").Value = 'athleta leggings with side pockets'
$ws.Range("A25").Value = 'athleta leggins'
$ws.Range("A26").Value = 'athleta mens'
$ws.Range("A27").Value = 'athleta pant'
$ws.Range("A28").Value = 'athleta pants'
$ws.Range("A29").Value = 'athleta pants women'
$ws.Range("A30").Value = 'athleta running pants'
$ws.Range("A31").Value = 'athleta running shorts'
$ws.Range("A32").Value = 'athleta shorts'
$ws.Range("A33").Value = 'athleta shorts with pockets women'
$ws.Range("A34").Value = 'athleta shorts women'
$ws.Range("A35").Value = 'athleta skirt'
$ws.Range("A36").Value = 'athleta tennis'
$ws.Range("A37").Value = 'athleta tennis skirt'
$ws.Range("A38").Value = 'athleta tights running'
$ws.Range("A39").Value = 'athleta tops women'
$ws.Range("A40").Value = 'athleta women'
$ws.Range("A41").Value = 'athlete leggings'
$ws.Range("A42").Value = 'athlete tape for pain'
$ws.Range("A43").Value = 'athletic apparel women'
$ws.Range("A44").Value = 'athletic black leggings'
$ws.Range("A45").Value = 'athletic capri'
$ws.Range("A46").Value = 'athletic capris'
$ws.Range("A47").Value = 'athletic capris for women'
$ws.Range("A48").Value = 'athletic capris for women with pockets'
$ws.Range("A49").Value = 'athletic cloth tape bulk'
$ws.Range("A50").Value = 'athletic clothes'
$ws.Range("A51").Value = 'athletic clothes for women'
$ws.Range("A52").Value = 'athletic clothes women'
$ws.Range("A53").Value = 'athletic cold weather leggings'
$ws.Range("A54").Value = 'athletic compression leggings'
$ws.Range("A55").Value = 'athletic compression leggings women'
$ws.Range("A56").Value = 'athletic compression pants'
$ws.Range("A57").Value = 'athletic compression shirt women'
$ws.Range("A58").Value = 'athletic compression tape'
$ws.Range("A59").Value = 'athletic compression tights'
$ws.Range("A60").Value = 'athletic compression wear'
$ws.Range("A61").Value = 'athletic gear for women'
$ws.Range("A62").Value = 'athletic high waist leggings'
$ws.Range("A63").Value = 'athletic knee compression'
$ws.Range("A64").Value = 'athletic knee support'
$ws.Range("A65").Value = 'athletic leggings capri'
$ws.Range("A66").Value = 'athletic leggings for women'
$ws.Range("A67").Value = 'athletic leggings for women capri'
$ws.Range("A68").Value = 'athletic leggings for women high waist'
$ws.Range("A69").Value = 'athletic leggings women'
$ws.Range("A70").Value = 'athletic leggings women black'
$ws.Range("A71").Value = 'athletic legings'
$ws.Range("A72").Value = 'athletic pants'
$ws.Range("A73").Value = 'athletic pants capri'
$ws.Range("A74").Value = 'athletic pants women'
$ws.Range("A75").Value = 'athletic running pants'
$ws.Range("A76").Value = 'athletic running pants women'
$ws.Range("A77").Value = 'athletic skirt with leggings'
$ws.Range("A78").Value = 'athletic skirts for women with leggings'
$ws.Range("A79").Value = 'athletic support'
$ws.Range("A80").Value = 'athletic support tape'
$ws.Range("A81").Value = 'athletic supports'
$ws.Range("A82").Value = 'athletic tape knee'
$ws.Range("A83").Value = 'athletic tape rainbow'
$ws.Range("A84").Value = 'athletic tape soccer'
$ws.Range("A85").Value = 'athletic tape weightlifting'
$ws.Range("A86").Value = 'athletic tight shorts women'
$ws.Range("A87").Value = 'athletic tights'
$ws.Range("A88").Value = 'athletic tights for women'
$ws.Range("A89").Value = 'athletic tights women'
$ws.Range("A90").Value = 'athletic training tape'
$ws.Range("A91").Value = 'athletic underwear women'
$ws.Range("A92").Value = 'athletic wear'
$ws.Range("A93").Value = 'athletic wear for women'
$ws.Range("A94").Value = 'athletic works capri'
$ws.Range("A95").Value = 'athletic works capri pants for women'
$ws.Range("A96").Value = 'athletic works capris for women'
$ws.Range("A97").Value = 'athletic yoga'
$ws.Range("A98").Value = 'athletics gear'
$ws.Range("A99").Value = 'athletics leggings'
$ws.Range("A100").Value = 'auto immune'
